$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the analysis timestamp (shared string referenced by A2)
$ws.Range("A2").Value = "2025-06-02 11:40:50"

# Update the metrics row (row 2) with the refreshed figures
$ws.Range("E2").Value = 212
$ws.Range("F2").Value = 8.313725490196077
$ws.Range("G2").Value = 1716
$ws.Range("H2").Value = 67.29411764705883
$ws.Range("I2").Value = 813
$ws.Range("J2").Value = 31.88235294117647
$ws.Range("K2").Value = 349184.6900000001
$ws.Range("N2").Value = 931.71
$ws.Range("O2").Value = 136
$ws.Range("P2").Value = 5.333333333333334
$ws.Range("Q2").Value = 13954.56
$ws.Range("R2").Value = 593
$ws.Range("S2").Value = 23.25490196078431
$ws.Range("T2").Value = 614
$ws.Range("U2").Value = 24.07843137254902
$ws.Range("V2").Value = 325460.06
$ws.Range("W2").Value = 305
$ws.Range("X2").Value = 11.96078431372549
$ws.Range("Y2").Value = 63
$ws.Range("Z2").Value = 2.470588235294117
$ws.Range("AA2").Value = 9770.07
$ws.Range("AB2").Value = 839
$ws.Range("AC2").Value = 32.90196078431373
$ws.Range("AE2").Value = 2513
$ws.Range("AF2").Value = 98.54901960784314
$ws.Range("AG2").Value = 37
$ws.Range("AH2").Value = 1.450980392156865
$ws.Range("AI2").Value = 59
$ws.Range("AJ2").Value = 81
$ws.Range("AK2").Value = 117
$ws.Range("AL2").Value = 22.95719844357977
$ws.Range("AM2").Value = 31.51750972762646
$ws.Range("AN2").Value = 45.52529182879378
$ws.Range("AO2").Value = 52656.29
$ws.Range("AP2").Value = 10118.19
$ws.Range("AQ2").Value = 3333.83
$ws.Range("AR2").Value = 79.65154456376212
$ws.Range("AS2").Value = 15.30547369914614
$ws.Range("AT2").Value = 5.042981737091752
$ws.Range("AU2").Value = 342.0262709802968
$ws.Range("AV2").Value = 814.1362916006339
$ws.Range("AW2").Value = 1332.513721185511
